# Arreglo algo que esta mal redactado
#
# Two wording fixes:
#  1) "... es una de tarea complicada y fundamental" -> "... es una tarea complicada y fundamental"
#     (drop the stray "de "), which in the canonical XML ends up as the original
#     run being split into three runs with identical rPr.
#  2) "lanzamientos que puede realizar un pitcher, a través de ..." gets split into
#     three runs (no text change) around the comma.
#
# Word's COM model merges/re-flows adjacent runs that share identical formatting
# whenever the underlying text of a paragraph is edited (deleted/replaced). To end
# up with the exact run layout from the target revision we:
#   1. make the actual text edit (deleting the stray "de"),
#   2. then re-impose every original run boundary in the affected paragraph
#      (plus the two new boundaries) by toggling a character property (Bold)
#      on/off across each sub-range -- this forces Word to split the run at
#      that boundary without leaving any residual formatting behind.

$d = $word.ActiveDocument

function Split-RunBoundary($doc, $base, $a, $b) {
    $rr = $doc.Range($base + $a, $base + $b)
    $rr.Bold = 1
    $rr.Bold = 0
}

# ---------------------------------------------------------------------------
# Hunk 1: "... los lanzamientos de un pitcher es una de tarea complicada ..."
# ---------------------------------------------------------------------------

$oldPhrase1 = " los lanzamientos de un pitcher es una de tarea complicada y fundamental"

$r1 = $d.Content
$r1.Find.ClearFormatting()
$found1 = $r1.Find.Execute($oldPhrase1, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found1) {
    throw "Could not find target phrase for hunk 1"
}

# $r1.Start now sits at local offset 413 within its paragraph (the text "La
# lectura de" + this run start the sentence); use that as our coordinate base.
$base1 = $r1.Start - 413

# Remove the stray "de" (and the space before "tarea" stays) - paragraph-local
# offsets 451..454 == " de" inside the matched phrase.
$delRange1 = $d.Range($base1 + 451, $base1 + 454)
$delRange1.Delete()

# Re-impose every original run boundary of the paragraph (it gets flattened by
# the Delete above because all of these runs share identical rPr), plus the
# two new boundaries needed for the 3-way split of the edited sentence.
Split-RunBoundary $d $base1 248 400
Split-RunBoundary $d $base1 400 413
Split-RunBoundary $d $base1 413 422
Split-RunBoundary $d $base1 422 451
Split-RunBoundary $d $base1 451 482
Split-RunBoundary $d $base1 482 511
Split-RunBoundary $d $base1 511 583
Split-RunBoundary $d $base1 583 584

# ---------------------------------------------------------------------------
# Hunk 2: "lanzamientos que puede realizar un pitcher, a través de ..."
# ---------------------------------------------------------------------------

$oldPhrase2 = "lanzamientos que puede realizar un pitcher, a través de la detección de los movimientos "

$r2 = $d.Content
$r2.Find.ClearFormatting()
$found2 = $r2.Find.Execute($oldPhrase2, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found2) {
    throw "Could not find target phrase for hunk 2"
}

# $r2.Start sits at local offset 35 within its paragraph.
$base2 = $r2.Start - 35

# No characters are removed here -- only new run boundaries are introduced
# around the comma -- but editing the paragraph still requires restoring the
# pre-existing boundaries too, so touch the whole run list once more.
Split-RunBoundary $d $base2 0 12
Split-RunBoundary $d $base2 12 25
Split-RunBoundary $d $base2 25 34
Split-RunBoundary $d $base2 34 35
Split-RunBoundary $d $base2 35 77
Split-RunBoundary $d $base2 77 78
Split-RunBoundary $d $base2 78 123
Split-RunBoundary $d $base2 123 138
Split-RunBoundary $d $base2 138 203
Split-RunBoundary $d $base2 203 221
Split-RunBoundary $d $base2 221 266

Write-Host "Done."
